{"js": "const replacements = [\n  [\"2024-10-03 Thursday\", \"2024-10-04 Friday\"],\n  [\"15\u00d764=960\", \"96\u00d726=2496\"],\n  [\"84\u00d734=2856\", \"66\u00d739=2574\"],\n  [\"19\u00d776=1444\", \"31\u00d716=496\"],\n  [\"60\u00d791=5460\", \"81\u00d794=7614\"],\n  [\"47\u00d721=987\", \"77\u00d763=4851\"],\n  [\"70\u00d751=3570\", \"58\u00d773=4234\"],\n  [\"21\u00d776=1596\", \"63\u00d765=4095\"],\n  [\"66\u00d780=5280\", \"41\u00d795=3895\"],\n  [\"56\u00d791=5096\", \"65\u00d732=2080\"],\n  [\"25\u00d739=975\", \"31\u00d781=2511\"],\n  [\"81\u00d756=4536\", \"30\u00d780=2400\"],\n  [\"52\u00d778=4056\", \"15\u00d791=1365\"],\n  [\"87\u00d747=4089\", \"76\u00d745=3420\"],\n  [\"36\u00d763=2268\", \"64\u00d742=2688\"],\n  [\"50\u00d759=2950\", \"61\u00d751=3111\"],\n  [\"40\u00d721=840\", \"22\u00d759=1298\"],\n  [\"22\u00d752=1144\", \"42\u00d792=3864\"],\n  [\"23\u00d790=2070\", \"84\u00d755=4620\"],\n  [\"91\u00d716=1456\", \"47\u00d717=799\"],\n  [\"16\u00d756=896\", \"50\u00d792=4600\"],\n  [\"34\u00d789=3026\", \"29\u00d782=2378\"],\n  [\"60\u00d786=5160\", \"53\u00d763=3339\"],\n  [\"78\u00d777=6006\", \"50\u00d723=1150\"],\n  [\"71\u00d754=3834\", \"37\u00d797=3589\"],\n  [\"64\u00d736=2304\", \"17\u00d788=1496\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$pairs = @(\n    @(\"2024-10-03 Thursday\", \"2024-10-04 Friday\"),\n    @(\"15\u00d764=960\", \"96\u00d726=2496\"),\n    @(\"84\u00d734=2856\", \"66\u00d739=2574\"),\n    @(\"19\u00d776=1444\", \"31\u00d716=496\"),\n    @(\"60\u00d791=5460\", \"81\u00d794=7614\"),\n    @(\"47\u00d721=987\", \"77\u00d763=4851\"),\n    @(\"70\u00d751=3570\", \"58\u00d773=4234\"),\n    @(\"21\u00d776=1596\", \"63\u00d765=4095\"),\n    @(\"66\u00d780=5280\", \"41\u00d795=3895\"),\n    @(\"56\u00d791=5096\", \"65\u00d732=2080\"),\n    @(\"25\u00d739=975\", \"31\u00d781=2511\"),\n    @(\"81\u00d756=4536\", \"30\u00d780=2400\"),\n    @(\"52\u00d778=4056\", \"15\u00d791=1365\"),\n    @(\"87\u00d747=4089\", \"76\u00d745=3420\"),\n    @(\"36\u00d763=2268\", \"64\u00d742=2688\"),\n    @(\"50\u00d759=2950\", \"61\u00d751=3111\"),\n    @(\"40\u00d721=840\", \"22\u00d759=1298\"),\n    @(\"22\u00d752=1144\", \"42\u00d792=3864\"),\n    @(\"23\u00d790=2070\", \"84\u00d755=4620\"),\n    @(\"91\u00d716=1456\", \"47\u00d717=799\"),\n    @(\"16\u00d756=896\", \"50\u00d792=4600\"),\n    @(\"34\u00d789=3026\", \"29\u00d782=2378\"),\n    @(\"60\u00d786=5160\", \"53\u00d763=3339\"),\n    @(\"78\u00d777=6006\", \"50\u00d723=1150\"),\n    @(\"71\u00d754=3834\", \"37\u00d797=3589\"),\n    @(\"64\u00d736=2304\", \"17\u00d788=1496\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $result = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $result) {\n        throw \"Find/Replace failed for: $oldText\"\n    }\n}"}
